$d = $word.ActiveDocument

# Locate the three consecutive paragraphs involved in this edit:
#   "Proof of Financial Responsibility. ..."   (kept)
#   "License Suspension. ..."                  (removed in full)
#   "Other Conditions. TESTS."                 (removed, except its own
#                                                trailing/empty closing run)
$proofPara = $null
$licensePara = $null
$otherPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "Proof of Financial Responsibility") {
        $proofPara = $p
    } elseif ($t -match "License Suspension") {
        $licensePara = $p
    } elseif ($t -match "Other Conditions") {
        $otherPara = $p
        break
    }
}

# 1) Remove the "License Suspension." paragraph and the "Other Conditions."
#    paragraph entirely (their text, runs and paragraph marks) in one shot:
#    from the start of the License Suspension paragraph through the end
#    (paragraph mark included) of the Other Conditions paragraph. This
#    leaves the Proof of Financial Responsibility paragraph's own mark
#    intact, so that paragraph keeps its original formatting, and the
#    blank paragraph that used to follow "Other Conditions." is untouched.
$cutRange = $d.Range($licensePara.Range.Start, $otherPara.Range.End)
$cutRange.Delete()

# 2) Clean up the trailing manual line break that was left dangling at the
#    end of the (now merged) "Proof of Financial Responsibility." paragraph.
$tail = $proofPara.Range
$tail.SetRange($tail.End - 2, $tail.End - 1)
if ([int][char]($tail.Text[0]) -eq 11) {
    $tail.Delete()
}
